# Delete the "2023-08-15" report row (row 2). This shifts the
# "2023-08-16" row up to row 2 and removes the now-unused shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()

# Update the Party Order Total / Net Profit values for the remaining row
$ws.Range("B2").Value = 205
$ws.Range("E2").Value = 205
